# Fill in the previously-empty "WEB EDU" result table (rows 49-53, columns B-O)
# on sheet "Feuil2" with the iteration-count / time(µs) pairs computed by the
# PageRank code, as described by the commit "tous les fichiers fait avec pagerank".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil2")

# alpha = 0.5
$ws.Range("B49").Value = 7
$ws.Range("C49").Value = 6389942
$ws.Range("D49").Value = 10
$ws.Range("E49").Value = 8818346
$ws.Range("F49").Value = 13
$ws.Range("G49").Value = 11436280
$ws.Range("H49").Value = 16
$ws.Range("I49").Value = 14033506
$ws.Range("J49").Value = 19
$ws.Range("K49").Value = 16945143
$ws.Range("L49").Value = 22
$ws.Range("M49").Value = 20457355
$ws.Range("N49").Value = 25
$ws.Range("O49").Value = 21925912

# alpha = 0.7
$ws.Range("B50").Value = 12
$ws.Range("C50").Value = 10657564
$ws.Range("D50").Value = 17
$ws.Range("E50").Value = 16672719
$ws.Range("F50").Value = 23
$ws.Range("G50").Value = 20395170
$ws.Range("H50").Value = 29
$ws.Range("I50").Value = 25786252
$ws.Range("J50").Value = 36
$ws.Range("K50").Value = 32647767
$ws.Range("L50").Value = 42
$ws.Range("M50").Value = 36374790
$ws.Range("N50").Value = 48
$ws.Range("O50").Value = 41605552

# alpha = 0.85
$ws.Range("B51").Value = 23
$ws.Range("C51").Value = 19996063
$ws.Range("D51").Value = 35
$ws.Range("E51").Value = 30315705
$ws.Range("F51").Value = 48
$ws.Range("G51").Value = 41606010
$ws.Range("H51").Value = 61
$ws.Range("I51").Value = 52643154
$ws.Range("J51").Value = 74
$ws.Range("K51").Value = 64128632
$ws.Range("L51").Value = 88
$ws.Range("M51").Value = 81659969
$ws.Range("N51").Value = 102
$ws.Range("O51").Value = 87412849

# alpha = 0.9
$ws.Range("B52").Value = 32
$ws.Range("C52").Value = 27586440
$ws.Range("D52").Value = 51
$ws.Range("E52").Value = 43790240
$ws.Range("F52").Value = 71
$ws.Range("G52").Value = 60538652
$ws.Range("H52").Value = 92
$ws.Range("I52").Value = 78105307
$ws.Range("J52").Value = 112
$ws.Range("K52").Value = 95042430
$ws.Range("L52").Value = 133
$ws.Range("M52").Value = 112512266
$ws.Range("N52").Value = 154
$ws.Range("O52").Value = 129769915

# alpha = 0.99
$ws.Range("B53").Value = 222
$ws.Range("C53").Value = 193505360
$ws.Range("D53").Value = 428
$ws.Range("E53").Value = 381549444
$ws.Range("F53").Value = 647
$ws.Range("G53").Value = 543966523
$ws.Range("H53").Value = 871
$ws.Range("I53").Value = 730146694
$ws.Range("J53").Value = 1097
$ws.Range("K53").Value = 917326174
$ws.Range("L53").Value = 1324
$ws.Range("M53").Value = 1109483471
$ws.Range("N53").Value = 1551
$ws.Range("O53").Value = 1297167096

# Reflect the updated view state recorded by Excel for the two sheets: Feuil1
# scrolled so column K is the left-most visible column, and Feuil2 scrolled so
# H30 is the top-left visible cell with O53 (the last value just entered) selected.
$ws1 = $wb.Worksheets.Item("Feuil1")
$ws1.Activate()
$win1 = $excel.ActiveWindow
$win1.ScrollColumn = 11
$win1.ScrollRow = 1

$ws.Activate()
$win2 = $excel.ActiveWindow
$win2.ScrollRow = 30
$win2.ScrollColumn = 8
[void]$ws.Range("O53").Select()
